$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume snapshot (Price = column D, Volume(1h) = column E).
# Values that look numeric are entered with a leading apostrophe so Excel keeps
# storing them as text (matching the sheet's existing inline-string Price column)
# instead of auto-converting them to numbers.
$ws.Range('D2').Value = '62.213.95'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.423.29'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'562.26"
$ws.Range('E5').Value = '  +2.53%  '
$ws.Range('D6').Value = "'144.20"
$ws.Range('E6').Value = '  +3.89%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'0.533"
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').Value = '2.421.71'
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  +1.78%  '
$ws.Range('D13').Value = "'0.353"
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').Value = "'26.11"
$ws.Range('E14').Value = '  +4.24%  '
$ws.Range('E15').Value = '  +6.00%  '
$ws.Range('D16').Value = '2.853.78'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').Value = '62.071.56'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '2.418.38'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('E19').Value = '  +2.86%  '

# Rows 20/21 swap ranking order: BitcoinCash moves up to rank 18, Polkadot drops to rank 19.
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'324.61"
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.19"
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = "'6.78"
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = "'65.51"
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = "'8.99"
$ws.Range('E26').Value = '  +6.34%  '
$ws.Range('D27').Value = "'588.27"
$ws.Range('E27').Value = '  +16.88%  '
$ws.Range('D28').Value = '2.542.94'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  +5.93%  '
$ws.Range('D31').Value = "'8.29"
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('D32').Value = "'1.45"
$ws.Range('E32').Value = '  +6.14%  '
$ws.Range('D33').Value = "'0.149"
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('E35').Value = '  +2.60%  '
$ws.Range('D36').Value = "'5.72"
$ws.Range('E36').Value = '  +5.41%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = "'4.78"
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').Value = "'154.00"
$ws.Range('E40').Value = '  +4.92%  '
$ws.Range('D41').Value = "'18.72"
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('E42').Value = '  -4.01%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  +10.84%  '
$ws.Range('D45').Value = "'150.95"
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').Value = "'3.65"
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('D47').Value = "'0.0539"
$ws.Range('E47').Value = '  +3.63%  '
$ws.Range('D48').Value = "'20.36"
$ws.Range('E48').Value = '  +5.53%  '
$ws.Range('D49').Value = "'0.592"
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').Value = "'0.0925"
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('D51').Value = "'0.0229"
$ws.Range('E51').Value = '  +2.68%  '
